$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("o19ool3dm6atef139bo63uo0a22rtab1", "2018-11-19T20:13:14.718"),
    @("40757p2y353be0086m0ee80usdbk4ind", "2018-11-19T20:15:50.197"),
    @("275280505v0mz9ci93z27h4xm3ivsm25", "2018-11-19T20:16:20.130"),
    @("rbd2c6564le4q4n6wcsi7zxoz0y118l6", "2018-11-19T20:17:20.081"),
    @("70o2abua894youqdmuk65eg490dym958", "2018-11-19T20:18:20.151"),
    @("63ks5b8lp4s8vhc474fz0a0974zd6cf8", "2018-11-19T20:18:32.616"),
    @("w1cd9g127myv39vl8s23vh18hta109by", "2018-11-19T20:18:47.222"),
    @("57baf3t70ulav9ewj720jxq8ph7z52nx", "2018-11-19T20:32:45.937"),
    @("owlze9c12wwgj7dm4pv1k1q5cdwun731", "2018-11-19T20:33:21.407"),
    @("dx7o5ofi69n9342qzpvmt8ov8657n3a4", "2018-11-19T20:34:22.887"),
    @("qnc5t5vy203opzg5cisdto8im9ffeknl", "2018-11-19T20:35:23.551"),
    @("wx9x7bp022rj85d7il1kmtd94r18wh44", "2018-11-19T20:36:20.032"),
    @("68m5z5s3133xwfa5snyra465x5127i1m", "2018-11-19T20:36:32.417"),
    @("by6di4s1rxu745k3l1qv56tzujlq5rvu", "2018-11-19T20:36:46.043"),
    @("ccpnf14573521j234spirk69tw6m6h51", "2018-11-19T20:38:22.423"),
    @("jv504elatekk3gw14xmoig5h0le4g6hl", "2018-11-19T20:39:24.893"),
    @("nddujp95o01av5z7l1cis8a8dx6tede2", "2018-11-19T20:40:20.356"),
    @("j9u7ln3io8ej618da904kzex7gjji5o3", "2018-11-19T20:41:20.524"),
    @("5z1i0303ehum0wb9374r8c8w8u4p9skj", "2018-11-19T20:42:20.875"),
    @("oa8b19zp12v35wl4fo2nyhl00ihm7n2x", "2018-11-19T20:42:33.768"),
    @("iswkwc81yq4xwthzm6jjy5up520h8wuy", "2018-11-19T20:42:46.369"),
    @("p8lqzngnlo38bc7obtk23k0hgvp4p9qy", "2018-11-19T20:42:59.886"),
    @("mgzi97jt7348f7c9r9sgx514e72uig4c", "2018-11-19T20:49:09.321"),
    @("5ugc9za232i12a96n4qgmh89c5mkp0e9", "2018-11-19T20:50:19.238"),
    @("4hx3676b4p4avl6m6vjp88kj6v7jr3qd", "2018-11-19T20:51:20.751"),
    @("d0j4m30z949vrgfum5qa8m6167ogycyn", "2018-11-19T20:52:21.182"),
    @("k9s46g8y8iaw68c0s2738ylk22nz8a8g", "2018-11-19T20:53:23.136"),
    @("74zfrkh58k0vve66q3j4518yhw2kk2kg", "2018-11-19T20:54:24.587"),
    @("gjmfd9l882kjn7ta0x1lrm8wycsc0m4s", "2018-11-19T20:55:20.404"),
    @("at5o7j9f90947s82lf7dz6qo21ab90ea", "2018-11-19T20:55:33.867"),
    @("ffm96i1l7mwqvf6c5pd1tkuxlwsv87o1", "2018-11-19T20:55:47.245"),
    @("3g6a5m06e34ja4ca40lr65izragi8lx3", "2018-11-19T20:56:22.357"),
    @("1s0u08gwd6ic3rbih7esvx1ruac1xrsu", "2018-11-19T20:57:22.670"),
    @("pov3asg3qdcxqei010m0c37g2v379j2x", "2018-11-19T20:58:23.020"),
    @("luruv8gam6xl6l6z9ewa6vj16o4n64jg", "2018-11-19T20:59:23.361"),
    @("svx9w05w711hjj8pnt48veqc56b7715w", "2018-11-19T21:00:25.872"),
    @("hqqduj297e13x2xd7bla9773c4ybn8to", "2018-11-19T21:01:20.653"),
    @("8fy6ppc4149tpxq742fm1k3h1rsc36k9", "2018-11-19T21:02:20.605"),
    @("0n8i1dri1ng4a6yx1odxx6qcd3q6w180", "2018-11-19T21:03:21.791"),
    @("o3js9tc8b713e1gprcgz2kg4iilo11m0", "2018-11-19T21:03:35.856"),
    @("8xo605rl3g78l4iurg8qfh23yll2lkkn", "2018-11-19T21:03:49.329"),
    @("o10lk8upoqi5curie4j8eh9lzo83kctt", "2018-11-19T21:04:03.014"),
)

$startRow = 611
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
